$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 16.40548401124649
$ws.Cells.Item(2, 4).Value = 9.538705566555747
$ws.Cells.Item(2, 5).Value = 16.55699451043912
$ws.Cells.Item(2, 6).Value = 39.03928311656505
$ws.Cells.Item(2, 7).Value = 3.693558962500307
$ws.Cells.Item(2, 10).Value = 12.27419186508692
$ws.Cells.Item(2, 11).Value = 8.630920874445264
$ws.Cells.Item(2, 12).Value = 8.569136173283496
$ws.Cells.Item(2, 13).Value = 15.30107565736081
$ws.Cells.Item(2, 15).Value = 29.39577894818369

$ws.Cells.Item(3, 2).Value = 16.32726379764847
$ws.Cells.Item(3, 4).Value = 9.548868416054141
$ws.Cells.Item(3, 5).Value = 16.5954759316503
$ws.Cells.Item(3, 6).Value = 39.13684199616634
$ws.Cells.Item(3, 7).Value = 3.695488555372158
$ws.Cells.Item(3, 10).Value = 12.29555224733683
$ws.Cells.Item(3, 11).Value = 8.393629060912421
$ws.Cells.Item(3, 12).Value = 8.544984415231752
$ws.Cells.Item(3, 13).Value = 15.27920299330132
$ws.Cells.Item(3, 15).Value = 29.47883508144989

$ws.Cells.Item(4, 2).Value = 16.28185121235945
$ws.Cells.Item(4, 4).Value = 9.55610785182664
$ws.Cells.Item(4, 5).Value = 16.62050056544015
$ws.Cells.Item(4, 6).Value = 39.20414676674118
$ws.Cells.Item(4, 7).Value = 3.696737030991849
$ws.Cells.Item(4, 10).Value = 12.30935206768452
$ws.Cells.Item(4, 11).Value = 8.245556643442365
$ws.Cells.Item(4, 12).Value = 8.53101383804292
$ws.Cells.Item(4, 13).Value = 15.26762074355413
$ws.Cells.Item(4, 15).Value = 29.53509894568174

$ws.Cells.Item(5, 2).Value = 16.26401737635405
$ws.Cells.Item(5, 4).Value = 9.559309793725266
$ws.Cells.Item(5, 5).Value = 16.63105041385833
$ws.Cells.Item(5, 6).Value = 39.23343346072058
$ws.Cells.Item(5, 7).Value = 3.697261862088848
$ws.Cells.Item(5, 10).Value = 12.31514822368516
$ws.Cells.Item(5, 11).Value = 8.184708494572369
$ws.Cells.Item(5, 12).Value = 8.525540121431874
$ws.Cells.Item(5, 13).Value = 15.26336917818731
$ws.Cells.Item(5, 15).Value = 29.55935005049983

$ws.Cells.Item(6, 2).Value = 16.26109708896054
$ws.Cells.Item(6, 4).Value = 9.55985669697516
$ws.Cells.Item(6, 5).Value = 16.63282349955329
$ws.Cells.Item(6, 6).Value = 39.23840873355869
$ws.Cells.Item(6, 7).Value = 3.697349981775309
$ws.Cells.Item(6, 10).Value = 12.31612111146977
$ws.Cells.Item(6, 11).Value = 8.174576781383124
$ws.Cells.Item(6, 12).Value = 8.524644565130254
$ws.Cells.Item(6, 13).Value = 15.26269159597868
$ws.Cells.Item(6, 15).Value = 29.56345680121444

$ws.Cells.Item(7, 2).Value = 16.28160795877099
$ws.Cells.Item(7, 4).Value = 9.556150014092589
$ws.Cells.Item(7, 5).Value = 16.62064141758126
$ws.Cells.Item(7, 6).Value = 39.20453421148775
$ws.Cells.Item(7, 7).Value = 3.696744043926145
$ws.Cells.Item(7, 10).Value = 12.30942953707947
$ws.Cells.Item(7, 11).Value = 8.24473795627735
$ws.Cells.Item(7, 12).Value = 8.530939124737635
$ws.Cells.Item(7, 13).Value = 15.26756150465092
$ws.Cells.Item(7, 15).Value = 29.53542064910711

$ws.Cells.Item(8, 2).Value = 16.3779799633948
$ws.Cells.Item(8, 4).Value = 9.542002598760787
$ws.Cells.Item(8, 5).Value = 16.56997349220892
$ws.Cells.Item(8, 6).Value = 39.07138364769477
$ws.Cells.Item(8, 7).Value = 3.694211095676202
$ws.Cells.Item(8, 10).Value = 12.28141516965749
$ws.Cells.Item(8, 11).Value = 8.549647499752428
$ws.Cells.Item(8, 12).Value = 8.560632556610873
$ws.Cells.Item(8, 13).Value = 15.29315276082404
$ws.Cells.Item(8, 15).Value = 29.42332275263377

$ws.Cells.Item(9, 2).Value = 16.58703544049263
$ws.Cells.Item(9, 4).Value = 9.522166898356312
$ws.Cells.Item(9, 5).Value = 16.48165893117539
$ws.Cells.Item(9, 6).Value = 38.86910652957289
$ws.Cells.Item(9, 7).Value = 3.689747133261724
$ws.Cells.Item(9, 10).Value = 12.23188683789329
$ws.Cells.Item(9, 11).Value = 9.125083186820365
$ws.Cells.Item(9, 12).Value = 8.625516118092667
$ws.Cells.Item(9, 13).Value = 15.35782725736526
$ws.Cells.Item(9, 15).Value = 29.24534775895011

$ws.Cells.Item(10, 2).Value = 16.75190691289655
$ws.Cells.Item(10, 4).Value = 9.512383394823834
$ws.Cells.Item(10, 5).Value = 16.42345383669598
$ws.Cells.Item(10, 6).Value = 38.75647066983558
$ws.Cells.Item(10, 7).Value = 3.686771009921443
$ws.Cells.Item(10, 10).Value = 12.19876257234016
$ws.Cells.Item(10, 11).Value = 9.52967053049119
$ws.Cells.Item(10, 12).Value = 8.677022549509257
$ws.Cells.Item(10, 13).Value = 15.41393716152676
$ws.Cells.Item(10, 15).Value = 29.1401758529348

$ws.Cells.Item(11, 2).Value = 16.82914265268558
$ws.Cells.Item(11, 4).Value = 9.508965697045079
$ws.Cells.Item(11, 5).Value = 16.39841386412985
$ws.Cells.Item(11, 6).Value = 38.71306120670353
$ws.Cells.Item(11, 7).Value = 3.685482331637001
$ws.Cells.Item(11, 10).Value = 12.18439535951574
$ws.Cells.Item(11, 11).Value = 9.708953889305095
$ws.Cells.Item(11, 12).Value = 8.701237319337835
$ws.Cells.Item(11, 13).Value = 15.4412731619528
$ws.Cells.Item(11, 15).Value = 29.09790047351137

$ws.Cells.Item(12, 2).Value = 16.85869197284265
$ws.Cells.Item(12, 4).Value = 9.507819392345731
$ws.Cells.Item(12, 5).Value = 16.38913777803132
$ws.Cells.Item(12, 6).Value = 38.69775014618479
$ws.Cells.Item(12, 7).Value = 3.685003663627922
$ws.Cells.Item(12, 10).Value = 12.17905518002362
$ws.Cells.Item(12, 11).Value = 9.776097885275631
$ws.Cells.Item(12, 12).Value = 8.710515131969297
$ws.Cells.Item(12, 13).Value = 15.45187991423144
$ws.Cells.Item(12, 15).Value = 29.08269356620663

$ws.Cells.Item(13, 2).Value = 16.85231489057909
$ws.Cells.Item(13, 4).Value = 9.508059701728847
$ws.Cells.Item(13, 5).Value = 16.3911263987131
$ws.Cells.Item(13, 6).Value = 38.70099750792316
$ws.Cells.Item(13, 7).Value = 3.685106339272923
$ws.Cells.Item(13, 10).Value = 12.18020082637742
$ws.Cells.Item(13, 11).Value = 9.761671351817879
$ws.Cells.Item(13, 12).Value = 8.70851225060048
$ws.Cells.Item(13, 13).Value = 15.44958429634312
$ws.Cells.Item(13, 15).Value = 29.08593296468725

$ws.Cells.Item(14, 2).Value = 16.83156775769327
$ws.Cells.Item(14, 4).Value = 9.508868429323565
$ws.Cells.Item(14, 5).Value = 16.39764659049312
$ws.Cells.Item(14, 6).Value = 38.71177895989222
$ws.Cells.Item(14, 7).Value = 3.685442764675423
$ws.Cells.Item(14, 10).Value = 12.18395401072913
$ws.Cells.Item(14, 11).Value = 9.714493163017348
$ws.Cells.Item(14, 12).Value = 8.70199846857639
$ws.Cells.Item(14, 13).Value = 15.44214070215591
$ws.Cells.Item(14, 15).Value = 29.09663331506049

$ws.Cells.Item(15, 2).Value = 16.81889825772289
$ws.Cells.Item(15, 4).Value = 9.509383040782955
$ws.Cells.Item(15, 5).Value = 16.40166720248493
$ws.Cells.Item(15, 6).Value = 38.71852973507993
$ws.Cells.Item(15, 7).Value = 3.685650048185186
$ws.Cells.Item(15, 10).Value = 12.18626600244852
$ws.Cells.Item(15, 11).Value = 9.685496190455861
$ws.Cells.Item(15, 12).Value = 8.698022542821343
$ws.Cells.Item(15, 13).Value = 15.43761436061624
$ws.Cells.Item(15, 15).Value = 29.10329204868976

$ws.Cells.Item(16, 2).Value = 16.74690268664947
$ws.Cells.Item(16, 4).Value = 9.512627480810961
$ws.Cells.Item(16, 5).Value = 16.42511912620258
$ws.Cells.Item(16, 6).Value = 38.75946524373745
$ws.Cells.Item(16, 7).Value = 3.686856535660567
$ws.Cells.Item(16, 10).Value = 12.19971557390613
$ws.Cells.Item(16, 11).Value = 9.517852661436319
$ws.Cells.Item(16, 12).Value = 8.675455450595521
$ws.Cells.Item(16, 13).Value = 15.41218669321685
$ws.Cells.Item(16, 15).Value = 29.14305079386405

$ws.Cells.Item(17, 2).Value = 16.70329438442411
$ws.Cells.Item(17, 4).Value = 9.514881947957427
$ws.Cells.Item(17, 5).Value = 16.4398738463056
$ws.Cells.Item(17, 6).Value = 38.78658412748793
$ws.Cells.Item(17, 7).Value = 3.687613336989894
$ws.Cells.Item(17, 10).Value = 12.20814571343898
$ws.Cells.Item(17, 11).Value = 9.413742651028615
$ws.Cells.Item(17, 12).Value = 8.661808922316082
$ws.Cells.Item(17, 13).Value = 15.39704789254476
$ws.Cells.Item(17, 15).Value = 29.16886845310908

$ws.Cells.Item(18, 2).Value = 16.6784236024267
$ws.Cells.Item(18, 4).Value = 9.516275915095976
$ws.Cells.Item(18, 5).Value = 16.44849574944501
$ws.Cells.Item(18, 6).Value = 38.80291904661419
$ws.Cells.Item(18, 7).Value = 3.68805476639493
$ws.Cells.Item(18, 10).Value = 12.21306053232659
$ws.Cells.Item(18, 11).Value = 9.353416074871213
$ws.Cells.Item(18, 12).Value = 8.654033897666451
$ws.Cells.Item(18, 13).Value = 15.38851127458855
$ws.Cells.Item(18, 15).Value = 29.18424209447988

$ws.Cells.Item(19, 2).Value = 16.67003970652082
$ws.Cells.Item(19, 4).Value = 9.51676461093658
$ws.Cells.Item(19, 5).Value = 16.45143825199445
$ws.Cells.Item(19, 6).Value = 38.80857627969402
$ws.Cells.Item(19, 7).Value = 3.688205282289561
$ws.Cells.Item(19, 10).Value = 12.21473595950283
$ws.Cells.Item(19, 11).Value = 9.332916028841289
$ws.Cells.Item(19, 12).Value = 8.651414271342466
$ws.Cells.Item(19, 13).Value = 15.38565041759844
$ws.Cells.Item(19, 15).Value = 29.18953730522004

$ws.Cells.Item(20, 2).Value = 16.70791480710457
$ws.Cells.Item(20, 4).Value = 9.514631894744777
$ws.Cells.Item(20, 5).Value = 16.43828917545747
$ws.Cells.Item(20, 6).Value = 38.78362100374338
$ws.Cells.Item(20, 7).Value = 3.687532139341163
$ws.Cells.Item(20, 10).Value = 12.20724148061559
$ws.Cells.Item(20, 11).Value = 9.424871884933628
$ws.Cells.Item(20, 12).Value = 8.663253982428778
$ws.Cells.Item(20, 13).Value = 15.39864180193416
$ws.Cells.Item(20, 15).Value = 29.16606587733889

$ws.Cells.Item(21, 2).Value = 16.83765366026358
$ws.Cells.Item(21, 4).Value = 9.508626877463414
$ws.Cells.Item(21, 5).Value = 16.39572586704002
$ws.Cells.Item(21, 6).Value = 38.70858158600856
$ws.Cells.Item(21, 7).Value = 3.685343695622445
$ws.Cells.Item(21, 10).Value = 12.18284888905475
$ws.Cells.Item(21, 11).Value = 9.728371267433133
$ws.Cells.Item(21, 12).Value = 8.703908824610536
$ws.Cells.Item(21, 13).Value = 15.44432018314329
$ws.Cells.Item(21, 15).Value = 29.09346859139009

$ws.Cells.Item(22, 2).Value = 16.92419537419773
$ws.Cells.Item(22, 4).Value = 9.505564037216885
$ws.Cells.Item(22, 5).Value = 16.36910871552625
$ws.Cells.Item(22, 6).Value = 38.66610936567189
$ws.Cells.Item(22, 7).Value = 3.68396776216036
$ws.Cells.Item(22, 10).Value = 12.16749175723879
$ws.Cells.Item(22, 11).Value = 9.922348739497044
$ws.Cells.Item(22, 12).Value = 8.731107611997972
$ws.Cells.Item(22, 13).Value = 15.47565856261873
$ws.Cells.Item(22, 15).Value = 29.05069641167972

$ws.Cells.Item(23, 2).Value = 16.87785289252994
$ws.Cells.Item(23, 4).Value = 9.507120091717052
$ws.Cells.Item(23, 5).Value = 16.38320519588746
$ws.Cells.Item(23, 6).Value = 38.68817603028725
$ws.Cells.Item(23, 7).Value = 3.68469716651187
$ws.Cells.Item(23, 10).Value = 12.17563478479139
$ws.Cells.Item(23, 11).Value = 9.819238612915839
$ws.Cells.Item(23, 12).Value = 8.716535139657424
$ws.Cells.Item(23, 13).Value = 15.45879860576674
$ws.Cells.Item(23, 15).Value = 29.07309665531321

$ws.Cells.Item(24, 2).Value = 16.70582528889273
$ws.Cells.Item(24, 4).Value = 9.51474463899695
$ws.Cells.Item(24, 5).Value = 16.43900517144905
$ws.Cells.Item(24, 6).Value = 38.78495831484289
$ws.Cells.Item(24, 7).Value = 3.687568829064093
$ws.Cells.Item(24, 10).Value = 12.20765007174303
$ws.Cells.Item(24, 11).Value = 9.419841824804962
$ws.Cells.Item(24, 12).Value = 8.662600450664209
$ws.Cells.Item(24, 13).Value = 15.39792067528907
$ws.Cells.Item(24, 15).Value = 29.16733126874161

$ws.Cells.Item(25, 2).Value = 16.52842774684365
$ws.Cells.Item(25, 4).Value = 9.526689521572283
$ws.Cells.Item(25, 5).Value = 16.50437349423211
$ws.Cells.Item(25, 6).Value = 38.91751705812141
$ws.Cells.Item(25, 7).Value = 3.690901218597007
$ws.Cells.Item(25, 10).Value = 12.24471000235444
$ws.Cells.Item(25, 11).Value = 8.972305861602676
$ws.Cells.Item(25, 12).Value = 8.607273110186313
$ws.Cells.Item(25, 13).Value = 15.33880364555764
$ws.Cells.Item(25, 15).Value = 29.28900649017866

